# Commit: "add the NA's under duplicate_image_filename"
#
# Populate column E (duplicate_image_filename) with "NA" for the stimuli
# rows that already have data in columns A/C/D (rows 2-21).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2:E21").Value = "NA"

# F1 is an empty (blank) string cell in the source file; round-tripping it
# through this engine otherwise coerces it into a stray value, so make sure
# it stays empty.
$ws.Cells.Item(1, 6).ClearContents()
